$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 116.9
$ws.Range("I12").Value = 147.75
$ws.Range("J12").Value = 96.333336
$ws.Range("K12").Value = 147.75
$ws.Range("L12").Value = 96.333336
$ws.Range("M12").Value = 22.25
$ws.Range("N12").Value = -436.333336

$ws.Range("H81").Value = 19990
$ws.Range("J81").Value = 19990
$ws.Range("L81").Value = 19990
$ws.Range("N81").Value = -21986

$ws.Range("H84").Value = 19990
$ws.Range("J84").Value = 19990
$ws.Range("L84").Value = 59970
$ws.Range("N84").Value = -69954

$ws.Range("H116").Value = 2318.8667
$ws.Range("I116").Value = 1970
$ws.Range("J116").Value = 2667.7334
$ws.Range("K116").Value = 1970
$ws.Range("L116").Value = 2667.7334
$ws.Range("M116").Value = 1472
$ws.Range("N116").Value = -9551.733400000001

$ws.Range("H137").Value = 2083.6667
$ws.Range("I137").Value = 999
$ws.Range("J137").Value = 2300.6
$ws.Range("K137").Value = 2997
$ws.Range("L137").Value = 6901.799999999999
$ws.Range("M137").Value = -447
$ws.Range("N137").Value = -12001.8

$ws.Range("H141").Value = 4149.05
$ws.Range("I141").Value = 889.3077
$ws.Range("J141").Value = 10202.857
$ws.Range("K141").Value = 2667.9231
$ws.Range("L141").Value = 30608.571
$ws.Range("M141").Value = 2512.0769
$ws.Range("N141").Value = -40968.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5001297.5
$ws.Range("I61").Value = 5001297.5
$ws.Range("K61").Value = 5001297.5
$ws.Range("M61").Value = -5001085.5

$ws.Range("H76").Value = 35861.332
$ws.Range("J76").Value = 35861.332
$ws.Range("L76").Value = 35861.332
$ws.Range("N76").Value = -36537.332

$ws.Range("H79").Value = 35861.332
$ws.Range("J79").Value = 35861.332
$ws.Range("L79").Value = 35861.332
$ws.Range("N79").Value = -38201.332

$ws.Range("H121").Value = 24500
$ws.Range("J121").Value = 24500
$ws.Range("L121").Value = 24500
$ws.Range("N121").Value = -27994

$ws.Range("H136").Value = 5001297.5
$ws.Range("I136").Value = 5001297.5
$ws.Range("K136").Value = 15003892.5
$ws.Range("M136").Value = -15001342.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 896469.5600000001
$ws.Range("I86").Value = 2410.5334
$ws.Range("J86").Value = 2115641
$ws.Range("K86").Value = 2410.5334
$ws.Range("L86").Value = 2115641
$ws.Range("M86").Value = -1287.5334
$ws.Range("N86").Value = -2117887

$ws.Range("H89").Value = 896469.5600000001
$ws.Range("I89").Value = 2410.5334
$ws.Range("J89").Value = 2115641
$ws.Range("K89").Value = 12052.667
$ws.Range("L89").Value = 10578205
$ws.Range("M89").Value = -6436.666999999999
$ws.Range("N89").Value = -10589437

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40414.285
$ws.Range("J28").Value = 40414.285
$ws.Range("L28").Value = 40414.285
$ws.Range("N28").Value = -40904.285

$ws.Range("H31").Value = 1079.862
$ws.Range("I31").Value = 874.14813
$ws.Range("J31").Value = 3857
$ws.Range("K31").Value = 874.14813
$ws.Range("L31").Value = 3857
$ws.Range("M31").Value = -579.14813
$ws.Range("N31").Value = -4447

$ws.Range("H34").Value = 1079.862
$ws.Range("I34").Value = 874.14813
$ws.Range("J34").Value = 3857
$ws.Range("K34").Value = 874.14813
$ws.Range("L34").Value = 3857
$ws.Range("M34").Value = -672.14813
$ws.Range("N34").Value = -4261

$ws.Range("H58").Value = 31250658
$ws.Range("I58").Value = 47619660
$ws.Range("J58").Value = 744.8182
$ws.Range("K58").Value = 47619660
$ws.Range("L58").Value = 744.8182
$ws.Range("M58").Value = -47619457
$ws.Range("N58").Value = -1150.8182

$ws.Range("H94").Value = 7005.5557
$ws.Range("I94").Value = 20371.8
$ws.Range("J94").Value = 1864.6923
$ws.Range("K94").Value = 20371.8
$ws.Range("L94").Value = 1864.6923
$ws.Range("M94").Value = -19920.8
$ws.Range("N94").Value = -2766.6923

$ws.Range("H136").Value = 31250658
$ws.Range("I136").Value = 47619660
$ws.Range("J136").Value = 744.8182
$ws.Range("K136").Value = 142858980
$ws.Range("L136").Value = 2234.4546
$ws.Range("M136").Value = -142856430
$ws.Range("N136").Value = -7334.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 29917724
$ws.Range("J5").Value = 7167.4375
$ws.Range("L5").Value = 21502.3125
$ws.Range("N5").Value = -21726.3125

$ws.Range("H63").Value = 798
$ws.Range("I63").Value = 798
$ws.Range("K63").Value = 2394
$ws.Range("M63").Value = -1645

$ws.Range("H66").Value = 798
$ws.Range("I66").Value = 798
$ws.Range("K66").Value = 7182
$ws.Range("M66").Value = -3438

$ws.Range("H68").Value = 9781.727999999999
$ws.Range("I68").Value = 600
$ws.Range("J68").Value = 11822.111
$ws.Range("K68").Value = 1800
$ws.Range("L68").Value = 35466.333
$ws.Range("M68").Value = -989
$ws.Range("N68").Value = -37088.333

$ws.Range("H71").Value = 9781.727999999999
$ws.Range("I71").Value = 600
$ws.Range("J71").Value = 11822.111
$ws.Range("K71").Value = 5400
$ws.Range("L71").Value = 106398.999
$ws.Range("M71").Value = -1344
$ws.Range("N71").Value = -114510.999

$ws.Range("H82").Value = 81
$ws.Range("I82").Value = 81
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 243
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = 163
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 81
$ws.Range("I85").Value = 81
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 243
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 1161
$ws.Range("N85").ClearContents()

$ws.Range("H94").Value = 1350
$ws.Range("I94").Value = 1350
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 4050
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -3374
$ws.Range("N94").ClearContents()

$ws.Range("H107").Value = 8074.7334
$ws.Range("I107").Value = 517.6667
$ws.Range("J107").Value = 9964
$ws.Range("K107").Value = 1553.0001
$ws.Range("L107").Value = 29892
$ws.Range("M107").Value = 366.9999
$ws.Range("N107").Value = -33732

$ws.Range("H112").Value = 66672036
$ws.Range("J112").Value = 100006100
$ws.Range("L112").Value = 300018300
$ws.Range("N112").Value = -300020516

$ws.Range("H121").Value = 10753279
$ws.Range("I121").Value = 233.55556
$ws.Range("J121").Value = 15152252
$ws.Range("K121").Value = 700.66668
$ws.Range("L121").Value = 45456756
$ws.Range("M121").Value = 609.33332
$ws.Range("N121").Value = -45459376

$ws.Range("H132").Value = 6887.3335
$ws.Range("I132").Value = 753.4286
$ws.Range("J132").Value = 10790.728
$ws.Range("K132").Value = 6780.8574
$ws.Range("L132").Value = 97116.552
$ws.Range("M132").Value = -4250.8574
$ws.Range("N132").Value = -102176.552

$ws.Range("H135").Value = 29917724
$ws.Range("J135").Value = 7167.4375
$ws.Range("L135").Value = 64506.9375
$ws.Range("N135").Value = -69576.9375
